$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.468.95'
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").Value = '2.193.20'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.02'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '83.66'
$ws.Range("E6").Value = '  +11.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.60'
$ws.Range("E10").Value = '  +10.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0917'
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.13'
$ws.Range("E12").Value = '  +5.76%  '
$ws.Range("E13").Value = '  +3.06%  '
$ws.Range("D14").Value = '2.521.14'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.33'
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").Value = '2.185.97'
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").Value = '43.367.60'
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.72'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.91'
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.37'
$ws.Range("E22").Value = '  +12.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.04'
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.91'
$ws.Range("E24").Value = '  -5.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.61'
$ws.Range("E26").Value = '  +1.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.45'
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.19'
$ws.Range("E28").Value = '  +5.55%  '
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("E30").Value = '  +3.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.26'
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.36'
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("E33").Value = '  +5.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.33'
$ws.Range("E34").Value = '  +3.54%  '
$ws.Range("E35").Value = '  +2.14%  '
$ws.Range("E36").Value = '  +4.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("E37").Value = '  +8.44%  '
$ws.Range("E38").Value = '  +7.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.41'
$ws.Range("E39").Value = '  +2.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.83'
$ws.Range("E40").Value = '  +9.36%  '
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("E42").Value = '  +7.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.72'
$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.28'
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0974'
$ws.Range("E47").Value = '  +0.40%  '
$ws.Range("E48").Value = '  +5.07%  '
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.438'
$ws.Range("E50").Value = '  -4.40%  '
$ws.Range("E51").Value = '  +11.45%  '
